# Auto-generated script applying numeric corrections to the per-job-profit
# worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as produced by the
# scheduled market data refresh runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 9200
$ws.Range("J88").Value = 9200
$ws.Range("L88").Value = 9200
$ws.Range("N88").Value = -10012
$ws.Range("H91").Value = 9200
$ws.Range("J91").Value = 9200
$ws.Range("L91").Value = 9200
$ws.Range("N91").Value = -12008
$ws.Range("H135").Value = 1055.3636
$ws.Range("I135").Value = 1015.9
$ws.Range("J135").Value = 1450
$ws.Range("K135").Value = 9143.1
$ws.Range("L135").Value = 13050
$ws.Range("M135").Value = -6608.1
$ws.Range("N135").Value = -18120
$ws.Range("H137").Value = 9758.200000000001
$ws.Range("I137").Value = 12376
$ws.Range("J137").Value = 3650
$ws.Range("K137").Value = 37128
$ws.Range("L137").Value = 10950
$ws.Range("M137").Value = -34578
$ws.Range("N137").Value = -16050
$ws.Range("H138").Value = 3319.5535
$ws.Range("J138").Value = 3822.718
$ws.Range("L138").Value = 11468.154
$ws.Range("N138").Value = -21748.154

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3122
$ws.Range("I45").Value = 1381.375
$ws.Range("K45").Value = 1381.375
$ws.Range("M45").Value = -1004.375
$ws.Range("H61").Value = 5015.2607
$ws.Range("I61").Value = 5197.5293
$ws.Range("K61").Value = 5197.5293
$ws.Range("M61").Value = -4985.5293
$ws.Range("H74").Value = 3100.1304
$ws.Range("I74").Value = 3169.3125
$ws.Range("J74").Value = 2942
$ws.Range("K74").Value = 3169.3125
$ws.Range("L74").Value = 2942
$ws.Range("M74").Value = -2295.3125
$ws.Range("N74").Value = -4690
$ws.Range("H77").Value = 3100.1304
$ws.Range("I77").Value = 3169.3125
$ws.Range("J77").Value = 2942
$ws.Range("K77").Value = 15846.5625
$ws.Range("L77").Value = 14710
$ws.Range("M77").Value = -11478.5625
$ws.Range("N77").Value = -23446
$ws.Range("H122").Value = 3484.1
$ws.Range("I122").Value = 3421.8
$ws.Range("J122").Value = 3671
$ws.Range("K122").Value = 10265.4
$ws.Range("L122").Value = 11013
$ws.Range("M122").Value = -7815.400000000001
$ws.Range("N122").Value = -15913
$ws.Range("H136").Value = 5015.2607
$ws.Range("I136").Value = 5197.5293
$ws.Range("K136").Value = 15592.5879
$ws.Range("M136").Value = -13042.5879

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2106.3572
$ws.Range("I99").Value = 898.9
$ws.Range("K99").Value = 898.9
$ws.Range("M99").Value = 599.1
$ws.Range("H106").Value = 21394.2
$ws.Range("J106").Value = 21394.2
$ws.Range("L106").Value = 21394.2
$ws.Range("N106").Value = -23918.2
$ws.Range("H107").Value = 7732.3335
$ws.Range("I107").Value = 1496.25
$ws.Range("J107").Value = 10000
$ws.Range("K107").Value = 1496.25
$ws.Range("L107").Value = 10000
$ws.Range("M107").Value = 423.75
$ws.Range("N107").Value = -13840
$ws.Range("H134").Value = 2975.923
$ws.Range("I134").Value = 2946.4285
$ws.Range("K134").Value = 8839.2855
$ws.Range("M134").Value = -6304.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2883.6667
$ws.Range("I31").Value = 1621.5714
$ws.Range("J31").Value = 3988
$ws.Range("K31").Value = 1621.5714
$ws.Range("L31").Value = 3988
$ws.Range("M31").Value = -1326.5714
$ws.Range("N31").Value = -4578
$ws.Range("H34").Value = 2883.6667
$ws.Range("I34").Value = 1621.5714
$ws.Range("J34").Value = 3988
$ws.Range("K34").Value = 1621.5714
$ws.Range("L34").Value = 3988
$ws.Range("M34").Value = -1419.5714
$ws.Range("N34").Value = -4392
$ws.Range("H58").Value = 1741.5714
$ws.Range("I58").Value = 1741.5714
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1741.5714
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1538.5714
$ws.Range("N58").Value = ""
$ws.Range("H99").Value = 2991.6667
$ws.Range("I99").Value = 2612.5
$ws.Range("K99").Value = 2612.5
$ws.Range("M99").Value = -1114.5
$ws.Range("H126").Value = 2991.6667
$ws.Range("I126").Value = 2612.5
$ws.Range("K126").Value = 7837.5
$ws.Range("M126").Value = -5367.5
$ws.Range("H136").Value = 1741.5714
$ws.Range("I136").Value = 1741.5714
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5224.7142
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2674.7142
$ws.Range("N136").Value = ""
$ws.Range("H140").Value = 131285.14
$ws.Range("J140").Value = 131285.14
$ws.Range("L140").Value = 131285.14
$ws.Range("N140").Value = -141645.14
$ws.Range("H141").Value = 49922.77
$ws.Range("I141").Value = 30000
$ws.Range("J141").Value = 94749
$ws.Range("K141").Value = 30000
$ws.Range("L141").Value = 94749
$ws.Range("M141").Value = -24820
$ws.Range("N141").Value = -105109

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 267.6
$ws.Range("I64").Value = 256
$ws.Range("J64").Value = 314
$ws.Range("K64").Value = 768
$ws.Range("L64").Value = 942
$ws.Range("M64").Value = -498
$ws.Range("N64").Value = -1482
$ws.Range("H67").Value = 267.6
$ws.Range("I67").Value = 256
$ws.Range("J67").Value = 314
$ws.Range("K67").Value = 768
$ws.Range("L67").Value = 942
$ws.Range("M67").Value = 168
$ws.Range("N67").Value = -2814
$ws.Range("H68").Value = 1496.5769
$ws.Range("J68").Value = 1893.7142
$ws.Range("L68").Value = 5681.142599999999
$ws.Range("N68").Value = -7303.142599999999
$ws.Range("H71").Value = 1496.5769
$ws.Range("J71").Value = 1893.7142
$ws.Range("L71").Value = 17043.4278
$ws.Range("N71").Value = -25155.4278
$ws.Range("H92").Value = 4724.875
$ws.Range("I92").Value = 3659.8
$ws.Range("K92").Value = 10979.4
$ws.Range("M92").Value = -9731.400000000001
$ws.Range("H107").Value = 4181.75
$ws.Range("I107").Value = 4113.5
$ws.Range("J107").Value = 4250
$ws.Range("K107").Value = 12340.5
$ws.Range("L107").Value = 12750
$ws.Range("M107").Value = -10420.5
$ws.Range("N107").Value = -16590
$ws.Range("H123").Value = 10000
$ws.Range("I123").Value = 10000
$ws.Range("K123").Value = 30000
$ws.Range("M123").Value = -27550
$ws.Range("H133").Value = 8965
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 37044004
$ws.Range("I113").Value = 90912060
$ws.Range("K113").Value = 90912060
$ws.Range("M113").Value = -90909890
$ws.Range("H136").Value = 38442
$ws.Range("J136").Value = 38442
$ws.Range("L136").Value = 115326
$ws.Range("N136").Value = -120426

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3071.182
$ws.Range("I136").Value = 3069.8572
$ws.Range("K136").Value = 9209.571599999999
$ws.Range("M136").Value = -6659.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 9676.333000000001
$ws.Range("J45").Value = 9611.6
$ws.Range("L45").Value = 9611.6
$ws.Range("N45").Value = -10593.6
$ws.Range("H74").Value = 22319.334
$ws.Range("J74").Value = 22319.334
$ws.Range("L74").Value = 22319.334
$ws.Range("N74").Value = -24191.334
$ws.Range("H77").Value = 22319.334
$ws.Range("J77").Value = 22319.334
$ws.Range("L77").Value = 66958.00199999999
$ws.Range("N77").Value = -76318.00199999999
$ws.Range("H122").Value = 5731.5454
$ws.Range("I122").Value = 4474.5
$ws.Range("K122").Value = 13423.5
$ws.Range("M122").Value = -10973.5
$ws.Range("H132").Value = 1976.317
$ws.Range("I132").Value = 1846.8975
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 5540.6925
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -3010.6925
$ws.Range("N132").Value = -18560
$ws.Range("H136").Value = 1569.7241
$ws.Range("I136").Value = 1371.1852
$ws.Range("K136").Value = 4113.5556
$ws.Range("M136").Value = -1563.5556

